# cobt and sensor testcases added to jenkins pipeline
$wb = $excel.ActiveWorkbook

$wsAppControl = $wb.Worksheets.Item("AppControl")
$wsSmoke      = $wb.Worksheets.Item("smoke")

# --- smoke sheet: rename/renumber existing COBT test cases and add new sensor case ---
$wsSmoke.Range("A17").Value2 = "ZestIOT_2268_Verifying_COBT_For_DIALCelebi_User"
$wsSmoke.Range("A18").Value2 = "ZestIOT_2268_Verifying_COBT_For_GMR_HYD_AISATS_User"
$wsSmoke.Range("A19").Value2 = "ZestIOT_2268_Verifying_COBT_For_GMR_HYD_SG_User"
$wsSmoke.Range("A20").Value2 = "ZestIOT_2293_GMR_HYD_Sensor_And_Scheduled_data_Validation"

# new row 21 - additional sensor validation test case, run flag "Y" same as row above
$wsSmoke.Range("A21").Value2 = "ZestIOT_2294_GMR_HYD_SensorATD_And_Scheduled_data_Validation"
$wsSmoke.Range("B21").Value2 = "Y"
$wsSmoke.Range("B20").Copy()
$wsSmoke.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- AppControl sheet: set the email id used for the report and hyperlink it ---
$wsAppControl.Range("B25").Value2 = "stiyyagura@enhops.com"
$hlink = $wsAppControl.Hyperlinks.Add($wsAppControl.Range("B25"), "mailto:stiyyagura@enhops.com")
$wsAppControl.Range("B25").Style = "Hyperlink"

# --- view/selection state: smoke tab is now the active tab ---
$wsAppControl.Activate()
$wsAppControl.Range("A26").Select()

$wsSmoke.Activate()
$wsSmoke.Range("A18").Select()
